# Adds the two new claim rows (58 and 59) to the "NEW" sheet, as captured
# by the upstream automated export ("Actualizacion automatica del index.html
# y archivo Excel"). Columns A-L and O-P are plain text in the source data
# (even when the text looks like a number or a date, e.g. "-503", "5",
# "7/10/2025"), while M/N (Coordenada_X / Coordenada_Y) are real numbers.
#
# Plain `Range.Value = "..."` lets Excel's input-parser reinterpret
# number-looking / date-looking text (turning "7/10/2025" into a date
# serial, "-503" into a numeric -503, etc.), so the text columns are
# temporarily forced to the "@" (Text) number format before the values are
# typed in, guaranteeing they land as literal text. ClearFormats() afterwards
# drops that temporary formatting again (back to the sheet's default style)
# without touching the now-literal text values, so the new rows end up
# styled exactly like the rest of the sheet (no explicit style index).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$textCols1 = $ws.Range("A58:L59")
$textCols2 = $ws.Range("O58:P59")
$textCols1.NumberFormat = "@"
$textCols2.NumberFormat = "@"

# Row 58
$ws.Range("A58").Value = "-503"
$ws.Range("B58").Value = "7/10/2025"
$ws.Range("C58").Value = "Salguero 842"
$ws.Range("D58").Value = "5"
$ws.Range("E58").Value = "808148673"
$ws.Range("F58").Value = "NEW"
$ws.Range("G58").Value = "Pendiente"
$ws.Range("H58").Value = "Cambiar columna picada en la base"
$ws.Range("I58").Value = "1"
$ws.Range("J58").Value = "Cambio"
$ws.Range("K58").Value = "Sin equipos"
$ws.Range("L58").Value = "Pasante"
$ws.Range("M58").Value = -58.419166
$ws.Range("N58").Value = -34.600265
$ws.Range("O58").Value = "Almagro"
$ws.Range("P58").Value = "Capital Sur"

# Row 59
$ws.Range("A59").Value = "-504"
$ws.Range("B59").Value = "7/10/2025"
$ws.Range("C59").Value = "Ohiggins 1611"
$ws.Range("D59").Value = "13"
$ws.Range("E59").Value = "808148679"
$ws.Range("F59").Value = "NEW"
$ws.Range("G59").Value = "Pendiente"
$ws.Range("H59").Value = "Columna podrida en la base"
$ws.Range("I59").Value = "1"
$ws.Range("J59").Value = "Cambio"
$ws.Range("K59").Value = "Nodo Teco"
$ws.Range("L59").Value = "Pasante"
$ws.Range("M59").Value = -58.448993
$ws.Range("N59").Value = -34.564383
$ws.Range("O59").Value = "Colegiales"
$ws.Range("P59").Value = "Capital Norte"

$textCols1.ClearFormats()
$textCols2.ClearFormats()
